# Commit: "Fruta / hortaliza, semanal"
#
# Weekly update for "Fruta, Macroferia Regional de Talca - Palta":
# four new price records (week of 2021-09-10, serial 44449) are inserted
# at the top of the data block (new rows 481-484), pushing all the
# existing records down by four rows (old 481-528 -> new 485-532).
#
# The repeated / constant columns for this sheet (single market, region,
# product) are:
#   A = 5
#   B = "Macroferia Regional de Talca"
#   C = "Maule"
#   E = 7
#   F = "Fruta"
#   G = 100106
#   H = "Oleaginosos"
#   I = 100106002
#   J = "Palta"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the 4 new records: insert rows at 481-484, shifting
#    the existing 481-528 block down to 485-532.
$ws.Rows("481:484").Insert()

# 2. Populate the 4 newly inserted rows with the new weekly data.
$newRows = @(
    @{ Row = 481; D = 44449; K = "Hass"; L = "1a nueva(o)";  M = 150; N = 3000;  O = 3000;  P = 3000;  Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; S = 3000; T = 1 },
    @{ Row = 482; D = 44449; K = "Hass"; L = "2a nueva(o)";  M = 100; N = 2500;  O = 2500;  P = 2500;  Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; S = 2500; T = 1 },
    @{ Row = 483; D = 44449; K = "Hass"; L = "3a nueva (o)"; M = 60;  N = 2000;  O = 2000;  P = 2000;  Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; S = 2000; T = 1 },
    @{ Row = 484; D = 44449; K = "Hass"; L = "Especial";     M = 800; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos";            R = "Perú";    S = 2800; T = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = 5
    $ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value  = "Maule"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 7
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100106
    $ws.Cells.Item($row, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value  = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
